$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 229, shifting existing rows 229..358 down to 230..359.
$ws.Rows(229).Insert()

# Populate the newly inserted row 229 with the new record's data.
$ws.Cells.Item(229, 1).Value = 4
$ws.Cells.Item(229, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(229, 3).Value = "Los Lagos"
$ws.Cells.Item(229, 4).Value = 44830
$ws.Cells.Item(229, 5).Value = 10
$ws.Cells.Item(229, 6).Value = 100114014
$ws.Cells.Item(229, 7).Value = "Betarraga"
$ws.Cells.Item(229, 8).Value = "Sin especificar"
$ws.Cells.Item(229, 9).Value = "Primera"
$ws.Cells.Item(229, 10).Value = 500
$ws.Cells.Item(229, 11).Value = 1500
$ws.Cells.Item(229, 12).Value = 1500
$ws.Cells.Item(229, 13).Value = 1500
$ws.Cells.Item(229, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(229, 15).Value = "Región del Maule"
$ws.Cells.Item(229, 16).Value = 300
$ws.Cells.Item(229, 17).Value = 5
$ws.Cells.Item(229, 18).Value = "Hortaliza"
